# Generate Report for Handoff
#
# A new handoff has just completed for the 93168050-... file, in both the
# zh-cn and de-de target locales. Update the "Latest Handoff Datetime"
# column (D) on each locale's status row for that file with the new
# handoff timestamp.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D5").Value = "2016-03-01 07:16:36"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D5").Value = "2016-03-01 07:16:45"
